$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 48
$ws.Range("A48").Value = 111974134
$ws.Range("B48").Value = 90658
$ws.Range("D48").Value = "NT"
$ws.Range("E48").Value = 4361
$ws.Range("F48").Value = "Orange taggsvamp"
$ws.Range("G48").Value = "Hydnellum aurantiacum"
$ws.Range("H48").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("P48").Value = "Aloppmoarna i S, Jmt"
$ws.Range("Q48").Value = 439399.8222122483
$ws.Range("R48").Value = 6952207.441512506

# Row 49
$ws.Range("A49").Value = 111974133
$ws.Range("B49").Value = 90682
$ws.Range("D49").Value = "NT"
$ws.Range("E49").Value = 2059
$ws.Range("F49").Value = "Skrovlig taggsvamp"
$ws.Range("G49").Value = "Hydnellum scabrosum"
$ws.Range("H49").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("P49").Value = "Aloppmoarna i S, Jmt"
$ws.Range("Q49").Value = 439389.9449806474
$ws.Range("R49").Value = 6952220.480550999

# Row 50
$ws.Range("A50").Value = 111974029
$ws.Range("B50").Value = 88032
$ws.Range("D50").Value = "VU"
$ws.Range("E50").Value = 6276
$ws.Range("F50").Value = "Goliatmusseron"
$ws.Range("G50").Value = "Tricholoma matsutake"
$ws.Range("H50").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("P50").Value = "Aloppmoarna, Jmt"
$ws.Range("Q50").Value = 439334.7866423383
$ws.Range("R50").Value = 6952296.802153576

# Row 51
$ws.Range("A51").Value = 111974125
$ws.Range("B51").Value = 90660
$ws.Range("D51").Value = "NT"
$ws.Range("E51").Value = 4362
$ws.Range("F51").Value = "Blå taggsvamp"
$ws.Range("G51").Value = "Hydnellum caeruleum"
$ws.Range("H51").Value = "(Hornem.) P.Karst."
$ws.Range("P51").Value = "Aloppmoarna i S, Jmt"
$ws.Range("Q51").Value = 439278.8711310769
$ws.Range("R51").Value = 6952206.909989387

# Row 53
$ws.Range("A53").Value = 111974126
$ws.Range("B53").Value = 88032
$ws.Range("D53").Value = "VU"
$ws.Range("E53").Value = 6276
$ws.Range("F53").Value = "Goliatmusseron"
$ws.Range("G53").Value = "Tricholoma matsutake"
$ws.Range("H53").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("P53").Value = "Aloppmoarna i S, Jmt"
$ws.Range("Q53").Value = 439289.9461055733
$ws.Range("R53").Value = 6952209.002200785

